# Digital_Data_tables2.xlsx — "Add files via upload"
#
# The lookup tables on the Age / Language / Period / Channel sheets each have
# a two-column table ("Code" + an Arabic-label column). The Arabic-label
# column header used to be sheet-specific (Age_AR, Language_AR, Period_AR,
# Channel_AR) - same as on the Service sheet used to be, before it was
# renamed to the generic "ARABIC" that the Service sheet's table already
# uses. Bring the remaining four sheets in line with that same "ARABIC"
# header, and leave the selection sitting on the header cell (B1) the way it
# was left after typing the new header in. The previously active tab
# (Service) also gets left behind in favor of the Channel tab being the one
# shown/selected when the file is reopened.

$wb = $excel.ActiveWorkbook

$wsAge      = $wb.Worksheets.Item("Age")
$wsService  = $wb.Worksheets.Item("Service")
$wsLanguage = $wb.Worksheets.Item("Language")
$wsPeriod   = $wb.Worksheets.Item("Period")
$wsChannel  = $wb.Worksheets.Item("Channel")

# Rename the Arabic-label table column header (B1) on every sheet that still
# used the old per-sheet name. Editing the header cell text directly also
# renames the corresponding ListObject/table column automatically.
$wsAge.Range("B1").Value      = "ARABIC"
$wsLanguage.Range("B1").Value = "ARABIC"
$wsPeriod.Range("B1").Value   = "ARABIC"
$wsChannel.Range("B1").Value  = "ARABIC"

# Leave each sheet's selection on the header cell (B1), and make Channel the
# active/selected tab instead of Service.
$wsAge.Activate()
$wsAge.Range("B1").Select()

$wsService.Activate()
$wsService.Range("B1").Select()

$wsLanguage.Activate()
$wsLanguage.Range("B1").Select()

$wsPeriod.Activate()
$wsPeriod.Range("B1").Select()

$wsChannel.Activate()
$wsChannel.Range("B1").Select()
